$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D34").Value = 1559.62
$ws1.Range("I34").Value = 92.40000000000001
$ws1.Range("M45").Value = 785.08

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F34").Value = 1914.26
$ws2.Range("F45").Value = 785.08
$ws2.Range("F58").Value = 23931.39

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 3460.42
$ws3.Range("E3").Value = 1531.7632
$ws3.Range("F3").Value = 0.6931676706095241

$ws3.Range("D8").Value = 1092.3
$ws3.Range("E8").Value = -342.3
$ws3.Range("F8").Value = 1.4564

$ws3.Range("D16").Value = 16557.76
$ws3.Range("E16").Value = 23832.41
$ws3.Range("F16").Value = 0.40994529114386

$ws3.Range("D19").Value = 23992.97
$ws3.Range("E19").Value = 31416.73560036207
$ws3.Range("F19").Value = 0.4330102414376159
